{"js": "// Load all paragraphs in the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Remove the \"Meta description: ...\" paragraph that directly follows the\n//    \"Play Pharaoh's Reign Mini-Max for Free - Review\" H1 heading.\nconst metaParagraph = paragraphs.items.find((p) =>\n  p.text.indexOf(\"Meta description\") !== -1\n);\nif (metaParagraph) {\n  metaParagraph.delete();\n  await context.sync();\n}\n\n// Re-load paragraphs (the collection/indices shifted after the delete).\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 2) Locate the closing paragraph that holds the old AI image-generation\n//    prompt text (last paragraph of the document body) and, right before\n//    it, insert a new bold paragraph repeating the page title.\nconst promptParagraph = paragraphs.items.find((p) =>\n  p.text.indexOf(\"Please create an eye-catching feature image\") !== -1\n);\n\nconst titlePara = promptParagraph.insertParagraph(\n  \"Play Pharaoh's Reign Mini-Max for Free - Review\",\n  \"Before\"\n);\ntitlePara.font.bold = true;\ntitlePara.font.italic = false;\n\n// 3) Replace the old image-generation prompt text with the meta-description\n//    copy (the run keeps its existing italic formatting).\npromptParagraph.insertText(\n  \"Read our review of Pharaoh's Reign Mini-Max and play for free. Features, graphics, and winning potential of the game.\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the \"Meta description: ...\" paragraph that directly follows the\n#    \"Play Pharaoh's Reign Mini-Max for Free - Review\" H1 heading.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"Meta description*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 2) Find the closing paragraph that still holds the old AI image-generation\n#    prompt text (last paragraph of the document body), insert a new bold\n#    paragraph right before it repeating the page title, then replace the\n#    prompt paragraph's text with the meta-description copy.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"Please create an eye-catching feature image*\") {\n        $rng = $p.Range\n        $rng.InsertParagraphBefore()\n\n        $newPara = $d.Paragraphs.Item($i)\n        $newRng = $newPara.Range\n        $newRng.Text = \"Play Pharaoh's Reign Mini-Max for Free - Review\"\n        $titleLen = $newRng.End - $newRng.Start\n        $titleTextRange = $d.Range($newRng.Start, $newRng.Start + $titleLen - 1)\n        $titleTextRange.Font.Italic = $false\n        $titleTextRange.Font.Bold = $true\n\n        $promptPara = $d.Paragraphs.Item($i + 1)\n        $promptRng = $promptPara.Range\n        $promptLen = $promptRng.End - $promptRng.Start\n        $promptTextRange = $d.Range($promptRng.Start, $promptRng.Start + $promptLen - 1)\n        $promptTextRange.Text = \"Read our review of Pharaoh's Reign Mini-Max and play for free. Features, graphics, and winning potential of the game.\"\n        break\n    }\n}\n"}
